$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18; existing rows 18-31 shift down to 19-32.
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Range("A18").Value = 11
$ws.Range("B18").Value = "Vega Monumental Concepción"
$ws.Range("C18").Value = "Bíobío"
$ws.Range("D18").Value = 44719
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = 100114007
$ws.Range("G18").Value = "Jengibre"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 13000
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = 13400
$ws.Range("N18").Value = "$/caja 13 kilos"
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 1031
$ws.Range("Q18").Value = 13
$ws.Range("R18").Value = "Hortaliza"
